# Trade #10 closed at 2026-02-17 04:07:03 - unknown UNKNOWN +0.000%
#
# Appends the newly-closed trade #10 to the "All Trades" and "MarketMaking"
# sheets, and rolls the aggregate stats on "Summary" / "Strategy Status"
# forward to account for it (one more total trade, one more losing trade).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.03   # Current Capital
$summary.Range("B4").Value = 0.03      # Total P&L $
$summary.Range("B5").Value = 0.06      # Total P&L %
$summary.Range("B6").Value = 10        # Total Trades
$summary.Range("B8").Value = 4         # Losing Trades
$summary.Range("B9").Value = 40        # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.03     # Capital
$status.Range("D4").Value = 10         # Trades
$status.Range("E4").Value = 0.03       # P&L $
$status.Range("F4").Value = 0.03       # P&L %
$status.Range("G4").Value = 40         # Win Rate %

# ---------------------------------------------------------------------------
# All Trades + MarketMaking sheets - append trade #10 as row 11
# ---------------------------------------------------------------------------
$tradeRow = @(10, "2026-02-17", "04:06:58", "MarketMaking", "DOWN", 0.82, 0.8, `
    "CLOSED", -2.439, -0.02, 100.03, 0, 0, 0.6, `
    "Normal spread capture: 19600 bps", "early_exit", 0.11)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $tradeRow.Length; $i++) {
        $cell = $ws.Cells.Item(11, $i + 1)
        if ($i -eq 1) {
            # Column B holds a "yyyy-mm-dd" looking string (e.g. "2026-02-17").
            # Plain assignment lets Excel re-interpret it as a date serial +
            # date number format, same as typing it into the UI would. Force
            # it to stay literal text, matching every other row in the
            # column, then drop the Text number format the trick leaves
            # behind so the cell keeps the sheet's default (unstyled) look.
            $cell.NumberFormat = "@"
            $cell.Value = $tradeRow[$i]
            $cell.ClearFormats()
        } else {
            $cell.Value = $tradeRow[$i]
        }
    }
}
